# chore: update Sheets via scheduled runner
# Refreshes the market-price / profit figures (columns H:N) on a handful
# of leve rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR
# sheets of the Ragnarok_Profits workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4037.3
$ws.Range("I64").Value = 3956.3333
$ws.Range("K64").Value = 3956.3333
$ws.Range("M64").Value = -3708.3333
$ws.Range("H67").Value = 4037.3
$ws.Range("I67").Value = 3956.3333
$ws.Range("K67").Value = 3956.3333
$ws.Range("M67").Value = -3098.3333
$ws.Range("H107").Value = 869.5599999999999
$ws.Range("I107").Value = 562.3889
$ws.Range("J107").Value = 1659.4286
$ws.Range("K107").Value = 562.3889
$ws.Range("L107").Value = 1659.4286
$ws.Range("M107").Value = 1357.6111
$ws.Range("N107").Value = -5499.4286
$ws.Range("H113").Value = 8619.75
$ws.Range("I113").Value = 8619.75
$ws.Range("K113").Value = 8619.75
$ws.Range("M113").Value = -5365.75
$ws.Range("H124").Value = 99891.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 99891.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 99891.5
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -109711.5
$ws.Range("H125").Value = 4365.5
$ws.Range("J125").Value = 4365.5
$ws.Range("L125").Value = 39289.5
$ws.Range("N125").Value = -44209.5
$ws.Range("H137").Value = 2183.7334
$ws.Range("I137").Value = 2172
$ws.Range("K137").Value = 6516
$ws.Range("M137").Value = -3966

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 179.33333
$ws.Range("I5").Value = 110
$ws.Range("J5").Value = 318
$ws.Range("K5").Value = 110
$ws.Range("L5").Value = 318
$ws.Range("M5").Value = 2
$ws.Range("N5").Value = -542
$ws.Range("I61").Value = 26251770
$ws.Range("J61").Value = 6667999.5
$ws.Range("K61").Value = 26251770
$ws.Range("L61").Value = 6667999.5
$ws.Range("M61").Value = -26251558
$ws.Range("N61").Value = -6668423.5
$ws.Range("H132").Value = 2330511.5
$ws.Range("I132").Value = 4463.1714
$ws.Range("K132").Value = 13389.5142
$ws.Range("M132").Value = -10859.5142
$ws.Range("I136").Value = 26251770
$ws.Range("J136").Value = 6667999.5
$ws.Range("K136").Value = 78755310
$ws.Range("L136").Value = 20003998.5
$ws.Range("M136").Value = -78752760
$ws.Range("N136").Value = -20009098.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 179.33333
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 318
$ws.Range("K4").Value = 110
$ws.Range("L4").Value = 318
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = -548
$ws.Range("H86").Value = 61145.082
$ws.Range("I86").Value = 135326
$ws.Range("K86").Value = 135326
$ws.Range("M86").Value = -134203
$ws.Range("H89").Value = 61145.082
$ws.Range("I89").Value = 135326
$ws.Range("K89").Value = 676630
$ws.Range("M89").Value = -671014
$ws.Range("H94").Value = 2110.4092
$ws.Range("I94").Value = 2129.4736
$ws.Range("J94").Value = 1989.6666
$ws.Range("K94").Value = 2129.4736
$ws.Range("L94").Value = 1989.6666
$ws.Range("M94").Value = -1678.4736
$ws.Range("N94").Value = -2891.6666
$ws.Range("H107").Value = 2422
$ws.Range("I107").Value = 2333.3076
$ws.Range("K107").Value = 2333.3076
$ws.Range("M107").Value = -413.3076000000001
$ws.Range("H134").Value = 5558131.5
$ws.Range("J134").Value = 25002322
$ws.Range("L134").Value = 75006966
$ws.Range("N134").Value = -75012036

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 47623070
$ws.Range("I31").Value = 83337420
$ws.Range("J31").Value = 3930.4443
$ws.Range("K31").Value = 83337420
$ws.Range("L31").Value = 3930.4443
$ws.Range("M31").Value = -83337125
$ws.Range("N31").Value = -4520.4443
$ws.Range("H34").Value = 47623070
$ws.Range("I34").Value = 83337420
$ws.Range("J34").Value = 3930.4443
$ws.Range("K34").Value = 83337420
$ws.Range("L34").Value = 3930.4443
$ws.Range("M34").Value = -83337218
$ws.Range("N34").Value = -4334.4443
$ws.Range("H94").Value = 983.17645
$ws.Range("J94").Value = 1027
$ws.Range("L94").Value = 1027
$ws.Range("N94").Value = -1929
$ws.Range("H99").Value = 30991.834
$ws.Range("I99").Value = 9334.333000000001
$ws.Range("J99").Value = 52649.332
$ws.Range("K99").Value = 9334.333000000001
$ws.Range("L99").Value = 52649.332
$ws.Range("M99").Value = -7836.333000000001
$ws.Range("N99").Value = -55645.332
$ws.Range("H105").Value = 1969.8636
$ws.Range("I105").Value = 1415
$ws.Range("K105").Value = 1415
$ws.Range("M105").Value = 332
$ws.Range("H126").Value = 30991.834
$ws.Range("I126").Value = 9334.333000000001
$ws.Range("J126").Value = 52649.332
$ws.Range("K126").Value = 28002.999
$ws.Range("L126").Value = 157947.996
$ws.Range("M126").Value = -25532.999
$ws.Range("N126").Value = -162887.996
$ws.Range("H134").Value = 3699.8333
$ws.Range("I134").Value = 3699.8333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11099.4999
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8564.499899999999
$ws.Range("N134").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 6765.6665
$ws.Range("I137").Value = 1723.1
$ws.Range("K137").Value = 5169.299999999999
$ws.Range("M137").Value = -69.29999999999927

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 53.8
$ws.Range("I2").Value = 48
$ws.Range("K2").Value = 48
$ws.Range("M2").Value = 65
$ws.Range("H22").Value = 10000
$ws.Range("I22").Value = 15000
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 15000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -14471
$ws.Range("N22").Value = -6058
$ws.Range("H107").Value = 1311.9546
$ws.Range("J107").Value = 1399
$ws.Range("L107").Value = 1399
$ws.Range("N107").Value = -5239
$ws.Range("H126").Value = 4102.75
$ws.Range("I126").Value = 3414.0833
$ws.Range("K126").Value = 10242.2499
$ws.Range("M126").Value = -7772.249899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 798.6
$ws.Range("I22").Value = 623.25
$ws.Range("K22").Value = 623.25
$ws.Range("M22").Value = -328.25
$ws.Range("H27").Value = 798.6
$ws.Range("I27").Value = 623.25
$ws.Range("K27").Value = 623.25
$ws.Range("M27").Value = -516.25
$ws.Range("H46").Value = 1974.6666
$ws.Range("I46").Value = 1499.5
$ws.Range("K46").Value = 1499.5
$ws.Range("M46").Value = -1311.5
$ws.Range("H97").Value = 61331.332
$ws.Range("J97").Value = 61331.332
$ws.Range("L97").Value = 61331.332
$ws.Range("N97").Value = -63313.332

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()
$ws.Range("H132").Value = 270123.9
$ws.Range("I132").Value = 6837.9688
$ws.Range("J132").Value = 1674315.6
$ws.Range("K132").Value = 20513.9064
$ws.Range("L132").Value = 5022946.800000001
$ws.Range("M132").Value = -17983.9064
$ws.Range("N132").Value = -5028006.800000001
$ws.Range("H136").Value = 799984.7
$ws.Range("I136").Value = 33316.75
$ws.Range("K136").Value = 99950.25
$ws.Range("M136").Value = -97400.25
